# file dialog and json conversion implemented
#
# Applies to 160603_ExcelLibraryEditor.xlsx:
#  - Construction sheet: F4 becomes a formula pulling the material name
#    from OpaqueMaterial!A35 (driving a file-open style lookup/conversion)
#    instead of a hard-coded shared string, and the active-cell/active-tab
#    bookkeeping moves to the OpaqueMaterial sheet.
#  - OpaqueMaterial sheet: eight new material rows are appended (63 is a
#    blank format-only spacer copied down from row 62), and the sheet view
#    now freezes the header row and keeps a scrolled-down selection.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("OpaqueMaterial")
$ws2 = $wb.Worksheets.Item("Construction")

# --- Construction sheet: replace the literal "XPS Board" text in F4 with
#     a formula that looks the name up from the OpaqueMaterial library ---
$ws2.Activate()
$ws2.Range("F4").Formula = "=OpaqueMaterial!A35"
$ws2.Range("E3").Select()

# --- OpaqueMaterial sheet: add a blank spacer row (format only, copied
#     down from the last existing row) then the new material rows ---
$ws1.Activate()

$ws1.Range("B62:L62").Copy()
$ws1.Range("B63").PasteSpecial(-4122)

$ws1.Cells.Item(64, 1).Value = " GypsumFibreBoard"
$ws1.Cells.Item(64, 2).Value = " Boards"
$ws1.Cells.Item(64, 3).Value = 0.32
$ws1.Cells.Item(64, 4).Value = 1000
$ws1.Cells.Item(64, 5).Value = 1100
$ws1.Cells.Item(64, 6).Value = 0.9
$ws1.Cells.Item(64, 7).Value = 0.7
$ws1.Cells.Item(64, 8).Value = 0.7
$ws1.Cells.Item(64, 9).Value = 45
$ws1.Cells.Item(64, 10).Value = 1.86
$ws1.Cells.Item(64, 11).Value = 0
$ws1.Cells.Item(64, 12).Value = " [lambda  rho c: Saint-Gobain Rigips][LCA  ICE (0.31fos + 0.41bio-embodied cabon emissions)]"

$ws1.Cells.Item(65, 1).Value = " Cross Laminated Timber"
$ws1.Cells.Item(65, 2).Value = " Timber"
$ws1.Cells.Item(65, 3).Value = 0.13
$ws1.Cells.Item(65, 4).Value = 500
$ws1.Cells.Item(65, 5).Value = 1600
$ws1.Cells.Item(65, 6).Value = 0.9
$ws1.Cells.Item(65, 7).Value = 0.7
$ws1.Cells.Item(65, 8).Value = 0.7
$ws1.Cells.Item(65, 9).Value = 10
$ws1.Cells.Item(65, 10).Value = 0.71
$ws1.Cells.Item(65, 11).Value = 0
$ws1.Cells.Item(65, 12).Value = " [lambda  rho c: dataholz.com][LCA  ICE (0.31fos + 0.41bio-embodied cabon emissions)]"

$ws1.Cells.Item(66, 1).Value = " Plaster"
$ws1.Cells.Item(66, 2).Value = " Screed"
$ws1.Cells.Item(66, 3).Value = 1
$ws1.Cells.Item(66, 4).Value = 2000
$ws1.Cells.Item(66, 5).Value = 1130
$ws1.Cells.Item(66, 6).Value = 0.9
$ws1.Cells.Item(66, 7).Value = 0.6
$ws1.Cells.Item(66, 8).Value = 0.6
$ws1.Cells.Item(66, 9).Value = 1.33
$ws1.Cells.Item(66, 10).Value = 0.221
$ws1.Cells.Item(66, 11).Value = 0
$ws1.Cells.Item(66, 12).Value = " [lambda  rho c: dataholz.com] [ LCA ICE"

$ws1.Cells.Item(67, 1).Value = " Mineral Wool"
$ws1.Cells.Item(67, 2).Value = " Insulation"
$ws1.Cells.Item(67, 3).Value = 0.041
$ws1.Cells.Item(67, 4).Value = 155
$ws1.Cells.Item(67, 5).Value = 1130
$ws1.Cells.Item(67, 6).Value = 0.9
$ws1.Cells.Item(67, 7).Value = 0.6
$ws1.Cells.Item(67, 8).Value = 0.6
$ws1.Cells.Item(67, 9).Value = 0
$ws1.Cells.Item(67, 10).Value = 0
$ws1.Cells.Item(67, 11).Value = 0
$ws1.Cells.Item(67, 12).Value = " [lambda  rho c: dataholz.com]"

$ws1.Cells.Item(68, 1).Value = " XPS Board"
$ws1.Cells.Item(68, 2).Value = " Insulation"
$ws1.Cells.Item(68, 3).Value = 0.034
$ws1.Cells.Item(68, 4).Value = 35
$ws1.Cells.Item(68, 5).Value = 1400
$ws1.Cells.Item(68, 6).Value = 0.9
$ws1.Cells.Item(68, 7).Value = 0.6
$ws1.Cells.Item(68, 8).Value = 0.6
$ws1.Cells.Item(68, 9).Value = 87.4
$ws1.Cells.Item(68, 10).Value = 2.8
$ws1.Cells.Item(68, 11).Value = 0
$ws1.Cells.Item(68, 12).Value = " "

$ws1.Cells.Item(69, 1).Value = " Sand-Lime Brick"
$ws1.Cells.Item(69, 2).Value = " Masonry"
$ws1.Cells.Item(69, 3).Value = 0.56
$ws1.Cells.Item(69, 4).Value = 1200
$ws1.Cells.Item(69, 5).Value = 1000
$ws1.Cells.Item(69, 6).Value = 0.9
$ws1.Cells.Item(69, 7).Value = 0.6
$ws1.Cells.Item(69, 8).Value = 0.6
$ws1.Cells.Item(69, 9).Value = 0
$ws1.Cells.Item(69, 10).Value = 0
$ws1.Cells.Item(69, 11).Value = 0
$ws1.Cells.Item(69, 12).Value = " "

$ws1.Cells.Item(70, 1).Value = " Bonded chippings"
$ws1.Cells.Item(70, 2).Value = " Screed"
$ws1.Cells.Item(70, 3).Value = 0.7
$ws1.Cells.Item(70, 4).Value = 1800
$ws1.Cells.Item(70, 5).Value = 1000
$ws1.Cells.Item(70, 6).Value = 0.9
$ws1.Cells.Item(70, 7).Value = 0.6
$ws1.Cells.Item(70, 8).Value = 0.6
$ws1.Cells.Item(70, 9).Value = 0
$ws1.Cells.Item(70, 10).Value = 0
$ws1.Cells.Item(70, 11).Value = 0
$ws1.Cells.Item(70, 12).Value = "  "

$ws1.Cells.Item(71, 1).Value = " Impact sound insulation"
$ws1.Cells.Item(71, 2).Value = " Insulation"
$ws1.Cells.Item(71, 3).Value = 0.035
$ws1.Cells.Item(71, 4).Value = 120
$ws1.Cells.Item(71, 5).Value = 1030
$ws1.Cells.Item(71, 6).Value = 0.9
$ws1.Cells.Item(71, 7).Value = 0.6
$ws1.Cells.Item(71, 8).Value = 0.6
$ws1.Cells.Item(71, 9).Value = 0
$ws1.Cells.Item(71, 10).Value = 0
$ws1.Cells.Item(71, 11).Value = 0
$ws1.Cells.Item(71, 12).Value = " "

# --- Freeze the header row on OpaqueMaterial and leave the view scrolled
#     down to the newly added rows, with OpaqueMaterial as the active tab ---
$ws1.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws1.Range("W87").Select()
